# "Update countries & provincias Spain" - refresh the COVID-19 "Pais" sheet
# with the next data pull (05:35 -> 06:52). The source table is kept sorted
# by "Casos totales" (column B) descending, so as per-country totals change
# a handful of rows swap rank and the country names in column A shift down
# to their new sorted position while picking up the refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp
$ws.Range("A1").Value = 'Datos actualizados a 26 de Junio de 2020 a las 06:52'

# Pakistan overtakes Alemania/Turquia (rows 15-17): Pakistan gets fresh
# totals, Alemania and Turquia shift down a rank but keep their own numbers.
$ws.Range("A15").Value = 'Pakistan'
$ws.Range("B15").Value = 195745
$ws.Range("C15").Value = 2775
$ws.Range("D15").Value = 84168
$ws.Range("E15").Value = 107615
$ws.Range("G15").Value = 59
$ws.Range("H15").Value = 3962

$ws.Range("A16").Value = 'Alemania'
$ws.Range("B16").Value = 193785
$ws.Range("D16").Value = 176800
$ws.Range("E16").Value = 7973
$ws.Range("H16").Value = 9012

$ws.Range("A17").Value = 'Turquia'
$ws.Range("B17").Value = 193115
$ws.Range("D17").Value = 165706
$ws.Range("E17").Value = 22363
$ws.Range("H17").Value = 5046

# Kazajistan (row 54): figures refreshed in place, no rank change.
$ws.Range("B54").Value = 19750
$ws.Range("C54").Value = 465
$ws.Range("E54").Value = 7390

# Haiti (row 81): figures refreshed in place, no rank change.
$ws.Range("B81").Value = 5543
$ws.Range("C81").Value = 114
$ws.Range("E81").Value = 4935
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 96

# Kirguistan overtakes Luxemburgo/Hungria (rows 90-92).
$ws.Range("A90").Value = 'Kirguistan'
$ws.Range("B90").Value = 4204
$ws.Range("C90").Value = 250
$ws.Range("D90").Value = 2162
$ws.Range("E90").Value = 1999
$ws.Range("H90").Value = 43

$ws.Range("A91").Value = 'Luxemburgo'
$ws.Range("B91").Value = 4151
$ws.Range("D91").Value = 3968
$ws.Range("E91").Value = 73
$ws.Range("H91").Value = 110

$ws.Range("A92").Value = 'Hungria'
$ws.Range("B92").Value = 4123
$ws.Range("D92").Value = 2640
$ws.Range("E92").Value = 906
$ws.Range("H92").Value = 577

# Gambia (row 190): figures refreshed in place, no rank change.
$ws.Range("B190").Value = 43
$ws.Range("C190").Value = 1
$ws.Range("E190").Value = 15

# Fiyi/Dominica tie (rows 202-203) swap order; identical totals either way.
$ws.Range("A202").Value = 'Fiyi'
$ws.Range("A203").Value = 'Dominica'

# Groenlandia/Islas Malvinas tie (rows 208-209) swap order; identical totals.
$ws.Range("A208").Value = 'Groenlandia'
$ws.Range("A209").Value = 'Islas Malvinas'

# Seychelles/Montserrat tie (rows 211-212) swap order; recovered/deaths
# figures (D/H) swap along with the names, rest of the row stays identical.
$ws.Range("A211").Value = 'Seychelles'
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = 'Montserrat'
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
